$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 705
$ws.Range("I12").Value = 627.6667
$ws.Range("J12").Value = 821
$ws.Range("K12").Value = 627.6667
$ws.Range("L12").Value = 821
$ws.Range("M12").Value = -457.6667
$ws.Range("N12").Value = -1161

$ws.Range("H64").Value = 13125
$ws.Range("I64").Value = 5000
$ws.Range("K64").Value = 5000
$ws.Range("M64").Value = -4752

$ws.Range("H67").Value = 13125
$ws.Range("I67").Value = 5000
$ws.Range("K67").Value = 5000
$ws.Range("M67").Value = -4142

$ws.Range("H74").Value = 0
$ws.Range("I74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("M74").ClearContents()

$ws.Range("H77").Value = 0
$ws.Range("I77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("M77").ClearContents()

$ws.Range("H138").Value = 3818.1333
$ws.Range("I138").Value = 3446.875
$ws.Range("J138").Value = 4242.4287
$ws.Range("K138").Value = 10340.625
$ws.Range("L138").Value = 12727.2861
$ws.Range("M138").Value = -5200.625
$ws.Range("N138").Value = -23007.2861

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 4754.857
$ws.Range("I45").Value = 1642
$ws.Range("K45").Value = 1642
$ws.Range("M45").Value = -1265

$ws.Range("H122").Value = 3053.6155
$ws.Range("I122").Value = 3231.6667
$ws.Range("K122").Value = 9695.000100000001
$ws.Range("M122").Value = -7245.000100000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 1000
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 1000
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 1000
$ws.Range("M5").ClearContents()
$ws.Range("N5").Value = -1226

$ws.Range("H22").Value = 318.9
$ws.Range("I22").Value = 373.33334
$ws.Range("J22").Value = 237.25
$ws.Range("K22").Value = 373.33334
$ws.Range("L22").Value = 237.25
$ws.Range("M22").Value = -200.33334
$ws.Range("N22").Value = -583.25

$ws.Range("H82").Value = 6628.5
$ws.Range("I82").Value = 6628.5
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 6628.5
$ws.Range("L82").Value = 0
$ws.Range("M82").Value = -6245.5
$ws.Range("N82").ClearContents()

$ws.Range("H85").Value = 6628.5
$ws.Range("I85").Value = 6628.5
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 6628.5
$ws.Range("L85").Value = 0
$ws.Range("M85").Value = -5302.5
$ws.Range("N85").ClearContents()

$ws.Range("H94").Value = 4976.1665
$ws.Range("I94").Value = 4928.5
$ws.Range("J94").Value = 5000
$ws.Range("K94").Value = 4928.5
$ws.Range("L94").Value = 5000
$ws.Range("M94").Value = -4477.5
$ws.Range("N94").Value = -5902

$ws.Range("H134").Value = 5030.2856
$ws.Range("I134").Value = 5662.4
$ws.Range("J134").Value = 3450
$ws.Range("K134").Value = 16987.2
$ws.Range("L134").Value = 10350
$ws.Range("M134").Value = -14452.2
$ws.Range("N134").Value = -15420

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1886.5555
$ws.Range("J31").Value = 1499
$ws.Range("L31").Value = 1499
$ws.Range("N31").Value = -2089

$ws.Range("H34").Value = 1886.5555
$ws.Range("J34").Value = 1499
$ws.Range("L34").Value = 1499
$ws.Range("N34").Value = -1903

$ws.Range("H41").Value = 28428.5
$ws.Range("J41").Value = 28428.5
$ws.Range("L41").Value = 28428.5
$ws.Range("N41").Value = -29284.5

$ws.Range("H58").Value = 1709.0769
$ws.Range("I58").Value = 1334.1
$ws.Range("K58").Value = 1334.1
$ws.Range("M58").Value = -1131.1

$ws.Range("H59").Value = 50000
$ws.Range("I59").Value = 0
$ws.Range("J59").Value = 50000
$ws.Range("K59").Value = 0
$ws.Range("L59").Value = 50000
$ws.Range("M59").ClearContents()
$ws.Range("N59").Value = -52290

$ws.Range("H62").Value = 2750

$ws.Range("H65").Value = 2750

$ws.Range("H68").Value = 5000
$ws.Range("I68").Value = 5000
$ws.Range("K68").Value = 5000
$ws.Range("M68").Value = -4251

$ws.Range("H70").Value = 24999.5
$ws.Range("J70").Value = 24999.5
$ws.Range("L70").Value = 24999.5
$ws.Range("N70").Value = -25629.5

$ws.Range("H71").Value = 5000
$ws.Range("I71").Value = 5000
$ws.Range("K71").Value = 15000
$ws.Range("M71").Value = -11256

$ws.Range("H73").Value = 24999.5
$ws.Range("J73").Value = 24999.5
$ws.Range("L73").Value = 24999.5
$ws.Range("N73").Value = -27183.5

$ws.Range("H74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("N74").ClearContents()

$ws.Range("H77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("N77").ClearContents()

$ws.Range("H136").Value = 1709.0769
$ws.Range("I136").Value = 1334.1
$ws.Range("K136").Value = 4002.3
$ws.Range("M136").Value = -1452.3

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 346.77777
$ws.Range("I8").Value = 346.77777
$ws.Range("K8").Value = 1040.33331
$ws.Range("M8").Value = -901.33331

$ws.Range("H38").Value = 273.375
$ws.Range("J38").Value = 309
$ws.Range("L38").Value = 927
$ws.Range("N38").Value = -1621

$ws.Range("H70").Value = 4475
$ws.Range("I70").Value = 950
$ws.Range("J70").Value = 8000
$ws.Range("K70").Value = 2850
$ws.Range("L70").Value = 24000
$ws.Range("M70").Value = -2535
$ws.Range("N70").Value = -24630

$ws.Range("H73").Value = 4475
$ws.Range("I73").Value = 950
$ws.Range("J73").Value = 8000
$ws.Range("K73").Value = 2850
$ws.Range("L73").Value = 24000
$ws.Range("M73").Value = -1758
$ws.Range("N73").Value = -26184

$ws.Range("H103").Value = 3999
$ws.Range("J103").Value = 3999
$ws.Range("L103").Value = 11997
$ws.Range("N103").Value = -13755

$ws.Range("H113").Value = 294.5
$ws.Range("J113").Value = 342.16666
$ws.Range("L113").Value = 1026.49998
$ws.Range("N113").Value = -5366.499980000001

$ws.Range("H131").Value = 1156.125
$ws.Range("I131").Value = 1099.5
$ws.Range("K131").Value = 3298.5
$ws.Range("M131").Value = 1741.5

$ws.Range("H134").Value = 166667330
$ws.Range("I134").Value = 166667330
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 500001990
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -499996920
$ws.Range("N134").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 0
$ws.Range("I80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("M80").ClearContents()

$ws.Range("H83").Value = 0
$ws.Range("I83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("M83").ClearContents()

$ws.Range("H102").Value = 1364.1666
$ws.Range("I102").Value = 1166.6666
$ws.Range("J102").Value = 1561.6666
$ws.Range("K102").Value = 1166.6666
$ws.Range("L102").Value = 1561.6666
$ws.Range("M102").Value = 455.3334
$ws.Range("N102").Value = -4805.6666

$ws.Range("H136").Value = 26217.334
$ws.Range("J136").Value = 26217.334
$ws.Range("L136").Value = 78652.00199999999
$ws.Range("N136").Value = -83752.00199999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2358.5
$ws.Range("I7").Value = 1949.625
$ws.Range("K7").Value = 1949.625
$ws.Range("M7").Value = -1837.625

$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("M22").ClearContents()

$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 0
$ws.Range("K27").Value = 0
$ws.Range("M27").ClearContents()

$ws.Range("H46").Value = 1836.44
$ws.Range("I46").Value = 1460.0625
$ws.Range("J46").Value = 2505.5557
$ws.Range("K46").Value = 1460.0625
$ws.Range("L46").Value = 2505.5557
$ws.Range("M46").Value = -1272.0625
$ws.Range("N46").Value = -2881.5557

$ws.Range("H55").Value = 434.69232
$ws.Range("I55").Value = 90.25
$ws.Range("J55").Value = 587.7778
$ws.Range("K55").Value = 90.25
$ws.Range("L55").Value = 587.7778
$ws.Range("M55").Value = 82.75
$ws.Range("N55").Value = -933.7778

$ws.Range("H126").Value = 2358.5
$ws.Range("I126").Value = 1949.625
$ws.Range("K126").Value = 5848.875
$ws.Range("M126").Value = -3378.875

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 9313
$ws.Range("I45").Value = 0
$ws.Range("J45").Value = 9313
$ws.Range("K45").Value = 0
$ws.Range("L45").Value = 9313
$ws.Range("M45").ClearContents()
$ws.Range("N45").Value = -10295

$ws.Range("H96").Value = 3214.1428
$ws.Range("J96").Value = 3566.6667
$ws.Range("L96").Value = 3566.6667
$ws.Range("N96").Value = -6312.6667
